$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the existing sheet and insert a new header row with "Counties" title
$ws.Name = "Counties"
$ws.Rows.Item(1).Insert()
$ws.Range("A1").Value = "Counties"
$ws.Range("A1").Font.Bold = $true
$null = $ws.Range("A2").Select()

# Add the new README sheet after Counties
$readme = $wb.Worksheets.Add([System.Type]::Missing, $ws)
$readme.Name = "README"

$readme.Range("A2").Value = "IEPR: 2024"
$readme.Range("A2").Font.Bold = $true
$readme.Range("A2").Characters(6, 5).Font.Bold = $false

$readme.Range("A4").Value = "Description:"
$readme.Range("A4").Font.Bold = $true

$readme.Range("A7").Value = "Usage:"
$readme.Range("A7").Font.Bold = $true

$readme.Range("A5").Value = "This describes the county level labels used in our data."

$readme.Range("A8").Value = "Double check any county-level data going into the model to make sure they match these labels."

$null = $readme.Range("A9").Select()
